$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.310.95'
$ws.Cells.Item(2, 5).Value = '  +0.58%  '
$ws.Cells.Item(3, 4).Value = '1.873.39'
$ws.Cells.Item(3, 5).Value = '  +0.52%  '
$ws.Cells.Item(4, 5).Value = '  +0.10%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.7115'
$ws.Cells.Item(5, 5).Value = '  +0.34%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '241.66'
$ws.Cells.Item(6, 5).Value = '  +0.10%  '
$ws.Cells.Item(7, 5).Value = '  +0.07%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3108'
$ws.Cells.Item(8, 5).Value = '  +0.54%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07772'
$ws.Cells.Item(9, 5).Value = '  +1.93%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '25.13'
$ws.Cells.Item(10, 5).Value = '  +2.13%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.08398'
$ws.Cells.Item(11, 5).Value = '  +0.48%  '
$ws.Cells.Item(12, 4).Value = '1.872.27'
$ws.Cells.Item(12, 5).Value = '  +0.30%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.240'
$ws.Cells.Item(13, 5).Value = '  +0.58%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.7113'
$ws.Cells.Item(14, 5).Value = '  +0.58%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '91.14'
$ws.Cells.Item(15, 5).Value = '  +0.02%  '
$ws.Cells.Item(16, 4).Value = '29.327.06'
$ws.Cells.Item(16, 5).Value = '  +0.53%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '6.060'
$ws.Cells.Item(17, 5).Value = '  +2.66%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000008196'
$ws.Cells.Item(18, 5).Value = '  +5.13%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '240.02'
$ws.Cells.Item(19, 5).Value = '  -0.98%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '13.20'
$ws.Cells.Item(20, 5).Value = '  +0.99%  '
$ws.Cells.Item(21, 4).Value = '2.119.10'
$ws.Cells.Item(21, 5).Value = '  +0.23%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '1.001'
$ws.Cells.Item(22, 5).Value = '  +0.16%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '7.756'
$ws.Cells.Item(23, 5).Value = '  -1.32%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '1.001'
$ws.Cells.Item(24, 5).Value = '  +0.13%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.1584'
$ws.Cells.Item(25, 5).Value = '  -0.12%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '163.07'
$ws.Cells.Item(26, 5).Value = '  -0.65%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.024'
$ws.Cells.Item(27, 5).Value = '  +0.90%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '18.49'
$ws.Cells.Item(28, 5).Value = '  +0.43%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.509'
$ws.Cells.Item(29, 5).Value = '  +0.67%  '
$ws.Cells.Item(30, 5).Value = '  +0.37%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.289'
$ws.Cells.Item(31, 5).Value = '  -2.49%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.319'
$ws.Cells.Item(32, 5).Value = '  +1.96%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.05294'
$ws.Cells.Item(33, 5).Value = '  +2.97%  '
$ws.Cells.Item(34, 5).Value = '  +1.61%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.177'
$ws.Cells.Item(35, 5).Value = '  +1.31%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.7441'
$ws.Cells.Item(36, 5).Value = '  -6.14%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.703'
$ws.Cells.Item(37, 5).Value = '  +0.62%  '
$ws.Cells.Item(38, 5).Value = '  +1.72%  '
$ws.Cells.Item(39, 4).Value = '1.226.67'
$ws.Cells.Item(39, 5).Value = '  +5.40%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.730'
$ws.Cells.Item(40, 5).Value = '  +1.28%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '6.478'
$ws.Cells.Item(41, 5).Value = '  +3.97%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.8854'
$ws.Cells.Item(42, 5).Value = '  -0.41%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '109.14'
$ws.Cells.Item(43, 5).Value = '  +6.36%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '72.46'
$ws.Cells.Item(44, 5).Value = '  -0.51%  '
$ws.Cells.Item(45, 5).Value = '  +0.09%  '
$ws.Cells.Item(46, 4).Value = '2.015.76'
$ws.Cells.Item(46, 5).Value = '  +0.33%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.798'
$ws.Cells.Item(47, 5).Value = '  +1.26%  '
$ws.Cells.Item(48, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.00000000123'
$ws.Cells.Item(48, 5).Value = '  +5.17%  '
$ws.Cells.Item(49, 2).Value = 'Mantle'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.5194'
$ws.Cells.Item(49, 5).Value = '  +0.25%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '9.381'
$ws.Cells.Item(50, 5).Value = '  +1.11%  '
$ws.Cells.Item(51, 5).Value = '  +0.94%  '
